$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in D17 with the new quiz "out of" points value
$ws.Range("D17").Value = 4

# F17 gets the Grade formula matching the pattern used in the rest of column F
$ws.Range("F17").Formula = "=C17/D17"

$ws.Activate()
$ws.Range("E4").Select()
